# Apply the "contingencies with rene fine" edit to the lines_states sheet.
#
# Two new line rows ("line7", "line8") are added right after "line6", and
# the contingency rows that used to follow ("extr1".."extr8") are updated
# in place / extended so the table grows from 15 to 17 data+header rows
# (A1:E15 -> A1:E17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target content for rows 8-17 (index, name, from_bus, to_bus,
# in_service). Rows 8-15 already exist and simply get new values written
# into them (no row-shifting needed); rows 16-17 are brand new and are
# appended below the current last row.
$data = @(
    @(8,  6,  "line7", 14, 11, $true),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $true),
    @(13, 11, "extr4", 7,  8,  $false),
    @(14, 12, "extr5", 9,  11, $true),
    @(15, 13, "extr6", 7,  11, $true),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $true)
)

foreach ($row in $data) {
    $r    = $row[0]
    $idx  = $row[1]
    $name = $row[2]
    $c    = $row[3]
    $d    = $row[4]
    $e    = $row[5]

    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}

# The new rows 16-17 need the same "index" column formatting (bold,
# bordered, centered) as the rest of column A. Copy the format from the
# already-styled A15 cell.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
